# The Dependency column (C) for rows that depend on the "AUTH" script
# was listing the free-text "Authentication" instead of the actual
# dependency name "AUTH" used in column A. Fix those cells so the
# dependency text matches the real row name.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($r in 3, 4, 5, 6, 8) {
    $ws.Cells.Item($r, 3).Value = "AUTH"
}
